# issue #5: stock data from json to db
# Update the stock ("股票") sheet to match the new normalized schema:
#  - drop the leading "★" marker from company names
#  - insert a new "category" column (value "normal" for every row)
#  - append "source_file" and "index" columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("B2").Value = "太平洋電線電纜股份有限公司"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-30"
$ws.Range("K2").Value = "黃志雄"
$ws.Range("L2").Value = 1366
$ws.Range("M2").Value = "tmpb8fa1"
$ws.Range("N2").Value = 72

$ws.Range("B3").Value = "華泰電子股份有限公司"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-03-30"
$ws.Range("K3").Value = "黃志雄"
$ws.Range("L3").Value = 1366
$ws.Range("M3").Value = "tmpb8fa1"
$ws.Range("N3").Value = 73

$ws.Range("B4").Value = "台灣光罩股份有限公司"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2012-03-30"
$ws.Range("K4").Value = "黃志雄"
$ws.Range("L4").Value = 1366
$ws.Range("M4").Value = "tmpb8fa1"
$ws.Range("N4").Value = 74

$ws.Range("B5").Value = "大同股份有限公司"
$ws.Range("I5").Value = "normal"
$ws.Range("J5").Value = "2012-03-30"
$ws.Range("K5").Value = "黃志雄"
$ws.Range("L5").Value = 1366
$ws.Range("M5").Value = "tmpb8fa1"
$ws.Range("N5").Value = 75

$ws.Range("B6").Value = "友達光電股份有限公司"
$ws.Range("I6").Value = "normal"
$ws.Range("J6").Value = "2012-03-30"
$ws.Range("K6").Value = "黃志雄"
$ws.Range("L6").Value = 1366
$ws.Range("M6").Value = "tmpb8fa1"
$ws.Range("N6").Value = 76

$ws.Range("B7").Value = "春雨開發股份有限公司(原正華）"
$ws.Range("I7").Value = "normal"
$ws.Range("J7").Value = "2012-03-30"
$ws.Range("K7").Value = "黃志雄"
$ws.Range("L7").Value = 1366
$ws.Range("M7").Value = "tmpb8fa1"
$ws.Range("N7").Value = 77

$ws.Range("B8").Value = "茂德科技股份有限公司"
$ws.Range("I8").Value = "normal"
$ws.Range("J8").Value = "2012-03-30"
$ws.Range("K8").Value = "黃志雄"
$ws.Range("L8").Value = 1366
$ws.Range("M8").Value = "tmpb8fa1"
$ws.Range("N8").Value = 78

$ws.Range("B9").Value = "金橋科技股份有限公司"
$ws.Range("I9").Value = "normal"
$ws.Range("J9").Value = "2012-03-30"
$ws.Range("K9").Value = "黃志雄"
$ws.Range("L9").Value = 1366
$ws.Range("M9").Value = "tmpb8fa1"
$ws.Range("N9").Value = 79

$ws.Range("B10").Value = "太平洋電線電纜股份有限公司"
$ws.Range("I10").Value = "normal"
$ws.Range("J10").Value = "2012-03-30"
$ws.Range("K10").Value = "黃志雄"
$ws.Range("L10").Value = 1366
$ws.Range("M10").Value = "tmpb8fa1"
$ws.Range("N10").Value = 80

$ws.Range("B11").Value = "中華映管股份有限公司"
$ws.Range("I11").Value = "normal"
$ws.Range("J11").Value = "2012-03-30"
$ws.Range("K11").Value = "黃志雄"
$ws.Range("L11").Value = 1366
$ws.Range("M11").Value = "tmpb8fa1"
$ws.Range("N11").Value = 81

$ws.Range("B12").Value = "台灣上地開發股份有限公司"
$ws.Range("I12").Value = "normal"
$ws.Range("J12").Value = "2012-03-30"
$ws.Range("K12").Value = "黃志雄"
$ws.Range("L12").Value = 1366
$ws.Range("M12").Value = "tmpb8fa1"
$ws.Range("N12").Value = 82

$ws.Range("B13").Value = "立端科技股份有限公司"
$ws.Range("I13").Value = "normal"
$ws.Range("J13").Value = "2012-03-30"
$ws.Range("K13").Value = "黃志雄"
$ws.Range("L13").Value = 1366
$ws.Range("M13").Value = "tmpb8fa1"
$ws.Range("N13").Value = 83

